$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 457. This shifts the existing rows
# 457-515 down to 459-517 (matching the data-shift pattern seen in the
# diff), and makes room for two brand-new weekly records at 457-458.
$ws.Rows("457:458").Insert()

# Row 457 - new "Primera" quality record dated 45131 (2023-07-24)
$ws.Cells.Item(457, 1).Value = 1
$ws.Cells.Item(457, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(457, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(457, 4).Value = 45131
$ws.Cells.Item(457, 5).Value = 15
$ws.Cells.Item(457, 6).Value = 100112043
$ws.Cells.Item(457, 7).Value = "Pepino ensalada"
$ws.Cells.Item(457, 8).Value = "Sin especificar"
$ws.Cells.Item(457, 9).Value = "Primera"
$ws.Cells.Item(457, 10).Value = 350
$ws.Cells.Item(457, 11).Value = 9000
$ws.Cells.Item(457, 12).Value = 10000
$ws.Cells.Item(457, 13).Value = 9714
$ws.Cells.Item(457, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(457, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(457, 16).Value = 139
$ws.Cells.Item(457, 17).Value = 70
$ws.Cells.Item(457, 18).Value = "Hortaliza"

# Row 458 - new "Segunda" quality record dated 45131 (2023-07-24)
$ws.Cells.Item(458, 1).Value = 1
$ws.Cells.Item(458, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(458, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(458, 4).Value = 45131
$ws.Cells.Item(458, 5).Value = 15
$ws.Cells.Item(458, 6).Value = 100112043
$ws.Cells.Item(458, 7).Value = "Pepino ensalada"
$ws.Cells.Item(458, 8).Value = "Sin especificar"
$ws.Cells.Item(458, 9).Value = "Segunda"
$ws.Cells.Item(458, 10).Value = 540
$ws.Cells.Item(458, 11).Value = 7000
$ws.Cells.Item(458, 12).Value = 8000
$ws.Cells.Item(458, 13).Value = 7556
$ws.Cells.Item(458, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(458, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(458, 16).Value = 76
$ws.Cells.Item(458, 17).Value = 100
$ws.Cells.Item(458, 18).Value = "Hortaliza"
